# Updated cryptos list on Wed Jun 14 17:09:08 UTC 2023 with GitHub Actions
#
# Refreshes the per-coin "Price" (column D) and "Volume(1h)" (column E) figures
# for rows 2-51 of Sheet1, and fixes the EnergySwap/Decentraland row order that
# had been swapped (rows 50-51), bringing each row Coin/Link/Price/Volume back
# in sync with the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings with European-style "." separators that look numeric;
# force text storage for the whole data range so Excel does not coerce them to floats.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.974.81'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '1.738.53'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '247.12'
$ws.Range("E5").Value = '  +4.29%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").Value = '0.5018'
$ws.Range("E7").Value = '  -2.07%  '

$ws.Range("D8").Value = '0.2727'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '0.06185'
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("D10").Value = '1.745.04'
$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '0.6534'
$ws.Range("E12").Value = '  +2.72%  '

$ws.Range("D13").Value = '15.16'
$ws.Range("E13").Value = '  +1.53%  '

$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '77.69'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("D18").Value = '25.987.17'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("D20").Value = '0.000006835'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").Value = '4.599'
$ws.Range("E21").Value = '  +8.11%  '

$ws.Range("D22").Value = '1.968.28'
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").Value = '8.770'
$ws.Range("E23").Value = '  +1.40%  '

$ws.Range("D24").Value = '5.410'
$ws.Range("E24").Value = '  +3.67%  '

$ws.Range("D25").Value = '134.00'
$ws.Range("E25").Value = '  -3.38%  '

$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").Value = '15.25'
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D28").Value = '1.787'
$ws.Range("E28").Value = '  +1.97%  '

$ws.Range("D29").Value = '105.22'
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").Value = '3.997'
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("D31").Value = '0.08141'
$ws.Range("E31").Value = '  -2.31%  '

$ws.Range("D32").Value = '3.697'
$ws.Range("E32").Value = '  +1.45%  '

$ws.Range("D33").Value = '0.04733'
$ws.Range("E33").Value = '  +3.85%  '

$ws.Range("D34").Value = '2.665'
$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").Value = '0.9961'
$ws.Range("E35").Value = '  +1.50%  '

$ws.Range("D36").Value = '0.6145'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").Value = '2.734'
$ws.Range("E37").Value = '  +1.68%  '

$ws.Range("D38").Value = '0.01615'
$ws.Range("E38").Value = '  +1.38%  '

$ws.Range("D39").Value = '1.942'
$ws.Range("E39").Value = '  +1.32%  '

$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").Value = '100.92'
$ws.Range("E41").Value = '  +3.21%  '

$ws.Range("D42").Value = '0.8069'
$ws.Range("E42").Value = '  +10.06%  '

$ws.Range("D43").Value = '0.3913'
$ws.Range("E43").Value = '  +2.05%  '

$ws.Range("E44").Value = '  +1.52%  '

$ws.Range("D45").Value = '0.1171'
$ws.Range("E45").Value = '  +4.16%  '

$ws.Range("D46").Value = '6.382'
$ws.Range("E46").Value = '  +3.57%  '

$ws.Range("D47").Value = '55.81'
$ws.Range("E47").Value = '  +2.02%  '

$ws.Range("D48").Value = '0.05291'
$ws.Range("E48").Value = '  +0.51%  '

$ws.Range("D49").Value = '30.90'
$ws.Range("E49").Value = '  +1.47%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.3480'
$ws.Range("E50").Value = '  +1.97%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.600'
$ws.Range("E51").Value = '  +0.90%  '

# Restore the default style on column D so only the cell content changed,
# matching the source edit (a data refresh, not a formatting change).
$ws.Range("D2:D51").Style = "Normal"
